$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for Price (D) column updates so numeric-looking
# strings (e.g. "412.28", "1.00") are not auto-converted to numbers by Excel.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

# Row 2
Set-TextValue $ws.Range("D2") "62.102.57"
$ws.Range("E2").Value = "  +0.16%  "

# Row 3
Set-TextValue $ws.Range("D3") "3.414.83"
$ws.Range("E3").Value = "  -0.14%  "

# Row 4
$ws.Range("E4").Value = "  -0.18%  "

# Row 5
Set-TextValue $ws.Range("D5") "412.28"
$ws.Range("E5").Value = "  +0.42%  "

# Row 6
Set-TextValue $ws.Range("D6") "128.67"
$ws.Range("E6").Value = "  -0.68%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.624"
$ws.Range("E7").Value = "  -2.80%  "

# Row 8
$ws.Range("E8").Value = "  +0.01%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.726"
$ws.Range("E9").Value = "  -1.80%  "

# Row 10
$ws.Range("E10").Value = "  -1.54%  "

# Row 11
Set-TextValue $ws.Range("D11") "42.70"
$ws.Range("E11").Value = "  -0.30%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.0000217"
$ws.Range("E12").Value = "  +0.21%  "

# Row 13
Set-TextValue $ws.Range("D13") "9.19"
$ws.Range("E13").Value = "  +0.70%  "

# Row 14
Set-TextValue $ws.Range("D14") "3.958.87"
$ws.Range("E14").Value = "  +0.10%  "

# Row 15
$ws.Range("E15").Value = "  -0.37%  "

# Row 16
Set-TextValue $ws.Range("D16") "20.48"
$ws.Range("E16").Value = "  -3.76%  "

# Row 17
Set-TextValue $ws.Range("D17") "3.418.22"
$ws.Range("E17").Value = "  +0.39%  "

# Row 18
Set-TextValue $ws.Range("D18") "12.74"
$ws.Range("E18").Value = "  +5.50%  "

# Row 19
$ws.Range("E19").Value = "  -1.24%  "

# Row 20
Set-TextValue $ws.Range("D20") "62.116.81"
$ws.Range("E20").Value = "  +0.21%  "

# Row 21
Set-TextValue $ws.Range("D21") "474.53"
$ws.Range("E21").Value = "  +7.18%  "

# Row 22
Set-TextValue $ws.Range("D22") "91.61"
$ws.Range("E22").Value = "  +0.56%  "

# Row 23
Set-TextValue $ws.Range("D23") "3.26"
$ws.Range("E23").Value = "  +2.44%  "

# Row 24
Set-TextValue $ws.Range("D24") "13.05"
$ws.Range("E24").Value = "  -0.65%  "

# Row 25
Set-TextValue $ws.Range("D25") "3.28"
$ws.Range("E25").Value = "  +0.45%  "

# Row 26
Set-TextValue $ws.Range("D26") "9.76"
$ws.Range("E26").Value = "  +10.28%  "

# Row 27
Set-TextValue $ws.Range("D27") "33.30"
$ws.Range("E27").Value = "  -0.86%  "

# Row 28
Set-TextValue $ws.Range("D28") "4.78"
$ws.Range("E28").Value = "  +0.56%  "

# Row 29
Set-TextValue $ws.Range("D29") "7.74"
$ws.Range("E29").Value = "  +1.59%  "

# Row 30
$ws.Range("E30").Value = "  -3.75%  "

# Row 31
Set-TextValue $ws.Range("D31") "11.83"
$ws.Range("E31").Value = "  -1.63%  "

# Row 32
$ws.Range("E32").Value = "  -1.44%  "

# Row 33
$ws.Range("E33").Value = "  -3.53%  "

# Row 34
Set-TextValue $ws.Range("D34") "40.94"
$ws.Range("E34").Value = "  -4.65%  "

# Row 35
Set-TextValue $ws.Range("D35") "1.00"
$ws.Range("E35").Value = "  +0.00%  "

# Row 36
Set-TextValue $ws.Range("D36") "58.05"
$ws.Range("E36").Value = "  +7.98%  "

# Row 37
Set-TextValue $ws.Range("D37") "0.0487"
$ws.Range("E37").Value = "  -3.39%  "

# Row 38
Set-TextValue $ws.Range("D38") "1.00"
$ws.Range("E38").Value = "  +0.15%  "

# Row 39
Set-TextValue $ws.Range("D39") "3.03"
$ws.Range("E39").Value = "  +3.74%  "

# Row 40
$ws.Range("E40").Value = "  -0.27%  "

# Row 41
Set-TextValue $ws.Range("D41") "0.322"
$ws.Range("E41").Value = "  +2.30%  "

# Row 42
Set-TextValue $ws.Range("D42") "146.65"
$ws.Range("E42").Value = "  +3.70%  "

# Row 43
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D43") "2.69"
$ws.Range("E43").Value = "  +11.96%  "

# Row 44
$ws.Range("B44").Value = "LidoDAOToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws.Range("D44") "3.31"
$ws.Range("E44").Value = "  -2.10%  "

# Row 45
Set-TextValue $ws.Range("D45") "2.06"
$ws.Range("E45").Value = "  +4.27%  "

# Row 46
Set-TextValue $ws.Range("D46") "4.32"
$ws.Range("E46").Value = "  +1.95%  "

# Row 47
Set-TextValue $ws.Range("D47") "2.33"
$ws.Range("E47").Value = "  +17.30%  "

# Row 48
Set-TextValue $ws.Range("D48") "16.31"
$ws.Range("E48").Value = "  -2.12%  "

# Row 49
Set-TextValue $ws.Range("D49") "0.0₃0530"
$ws.Range("E49").Value = "  +23.88%  "

# Row 50
Set-TextValue $ws.Range("D50") "22.23"
$ws.Range("E50").Value = "  +0.08%  "

# Row 51
Set-TextValue $ws.Range("D51") "113.69"
$ws.Range("E51").Value = "  +7.95%  "
